# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled
# update). Updates the Price (D) / Volume(1h) (E) columns for every coin
# row, and re-syncs two coin pairs whose relative ranking flipped this
# run (VeChain/Stellar at rows 40-41, Quant/Flow at rows 49-50) by
# swapping their Coin/Link/Price/Volume cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# remain text (matching the original sheet, where every Price cell is
# stored as a text string) - otherwise Excel auto-converts them to numbers.
$textForceCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the commit.
$ws.Range("D2").Value = '24.607.04'
$ws.Range("E2").Value = '  +3.17%  '
$ws.Range("D3").Value = '1.697.34'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '316.66'
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.3943'
$ws.Range("E7").Value = '  +1.52%  '
$ws.Range("D8").Value = '0.4020'
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("D9").Value = '1.534'
$ws.Range("E9").Value = '  +4.54%  '
$ws.Range("D10").Value = '0.9997'
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").Value = '53.04'
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("D12").Value = '0.08771'
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").Value = '7.229'
$ws.Range("E13").Value = '  +7.51%  '
$ws.Range("D14").Value = '23.29'
$ws.Range("E14").Value = '  +2.81%  '
$ws.Range("D15").Value = '8.056'
$ws.Range("E15").Value = '  +9.90%  '
$ws.Range("D16").Value = '0.00001316'
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("D17").Value = '1.694.59'
$ws.Range("E17").Value = '  +1.94%  '
$ws.Range("D18").Value = '99.93'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '0.07074'
$ws.Range("E19").Value = '  +2.83%  '
$ws.Range("D20").Value = '19.71'
$ws.Range("E20").Value = '  +3.26%  '
$ws.Range("D21").Value = '6.955'
$ws.Range("E21").Value = '  +4.72%  '
$ws.Range("D22").Value = '0.9991'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '14.20'
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("D24").Value = '24.597.13'
$ws.Range("E24").Value = '  +3.16%  '
$ws.Range("D25").Value = '3.170'
$ws.Range("E25").Value = '  +11.33%  '
$ws.Range("D26").Value = '2.338'
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").Value = '22.69'
$ws.Range("E27").Value = '  +4.40%  '
$ws.Range("D28").Value = '162.29'
$ws.Range("E28").Value = '  +2.12%  '
$ws.Range("D29").Value = '137.46'
$ws.Range("E29").Value = '  +5.65%  '
$ws.Range("D30").Value = '5.194'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("D31").Value = '7.572'
$ws.Range("E31").Value = '  +10.63%  '
$ws.Range("D32").Value = '1.881.76'
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("D33").Value = '1.095'
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("D34").Value = '0.08602'
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("D35").Value = '7.219'
$ws.Range("E35").Value = '  +10.11%  '
$ws.Range("D36").Value = '11.33'
$ws.Range("E36").Value = '  +8.18%  '
$ws.Range("D37").Value = '0.2739'
$ws.Range("E37").Value = '  +3.69%  '
$ws.Range("D38").Value = '1.930'
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.09095'
$ws.Range("E40").Value = '  +3.44%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.02746'
$ws.Range("E41").Value = '  +9.10%  '
$ws.Range("D42").Value = '1.477'
$ws.Range("E42").Value = '  +1.37%  '
$ws.Range("D43").Value = '0.7683'
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("D44").Value = '0.7199'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").Value = '15.67'
$ws.Range("E45").Value = '  +4.01%  '
$ws.Range("D46").Value = '2.561'
$ws.Range("E46").Value = '  +5.56%  '
$ws.Range("D47").Value = '4.215'
$ws.Range("E47").Value = '  +2.65%  '
$ws.Range("D48").Value = '0.9990'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D49").Value = '1.333'
$ws.Range("E49").Value = '  +8.92%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '141.13'
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("D51").Value = '0.08002'
$ws.Range("E51").Value = '  +2.84%  '

# Restore default (unstyled) formatting on the cells we force-formatted
# as text, so the workbook styling matches the original layout.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}